# Auto-generated Excel COM-interop script
# Applies cell value updates to the crypto-price worksheet, matching the
# upstream data-refresh commit. A handful of D-column prices are plain
# numeric-looking strings (e.g. "230.58"); the source data keeps these as
# literal text, so those specific cells are pre-formatted as Text before
# the value is written (otherwise Excel's normal typed-input parsing would
# turn them into real numbers). Every other cell is untouched stylistically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5,D7,D10,D11,D12,D13,D15,D16,D17,D21,D22,D23,D28,D29,D30,D35,D37,D39,D42,D43,D47,D48").NumberFormat = "@"

$ws.Range("D2").Value = '44.298.30'
$ws.Range("E2").Value = '  +5.17%  '
$ws.Range("D3").Value = '2.267.31'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '230.58'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  +2.72%  '
$ws.Range("D7").Value = '63.72'
$ws.Range("E7").Value = '  +4.65%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +11.08%  '
$ws.Range("D10").Value = '0.104'
$ws.Range("E10").Value = '  +15.87%  '
$ws.Range("D11").Value = '57.03'
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").Value = '26.34'
$ws.Range("E12").Value = '  +19.45%  '
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  +2.21%  '
$ws.Range("D14").Value = '2.603.74'
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("D15").Value = '15.74'
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Value = '6.13'
$ws.Range("E16").Value = '  +9.96%  '
$ws.Range("D17").Value = '0.838'
$ws.Range("E17").Value = '  +5.19%  '
$ws.Range("D18").Value = '2.270.34'
$ws.Range("E18").Value = '  +2.21%  '
$ws.Range("D19").Value = '44.112.50'
$ws.Range("E19").Value = '  +4.91%  '
$ws.Range("E20").Value = '  +9.12%  '
$ws.Range("D21").Value = '73.52'
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").Value = '6.03'
$ws.Range("E22").Value = '  -2.72%  '
$ws.Range("D23").Value = '251.60'
$ws.Range("E23").Value = '  +3.48%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("E26").Value = '  +2.18%  '
$ws.Range("E27").Value = '  +24.97%  '
$ws.Range("D28").Value = '10.00'
$ws.Range("E28").Value = '  +4.23%  '
$ws.Range("D29").Value = '172.15'
$ws.Range("E29").Value = '  +1.89%  '
$ws.Range("D30").Value = '20.80'
$ws.Range("E30").Value = '  +2.64%  '
$ws.Range("E31").Value = '  -2.01%  '
$ws.Range("E32").Value = '  -4.98%  '
$ws.Range("E33").Value = '  +3.29%  '
$ws.Range("E34").Value = '  +5.67%  '
$ws.Range("D35").Value = '4.74'
$ws.Range("E35").Value = '  +3.19%  '
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("D37").Value = '3.82'
$ws.Range("E37").Value = '  +7.63%  '
$ws.Range("E38").Value = '  +5.69%  '
$ws.Range("D39").Value = '2.31'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("E40").Value = '  +4.01%  '
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("D42").Value = '17.45'
$ws.Range("E42").Value = '  +8.25%  '
$ws.Range("D43").Value = '8.29'
$ws.Range("E43").Value = '  -2.91%  '
$ws.Range("E44").Value = '  +1.05%  '
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("B47").Value = 'TerraClassic'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D47").Value = '0.000210'
$ws.Range("E47").Value = '  -7.91%  '
$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").Value = '4.35'
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("D49").Value = '1.441.16'
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("E50").Value = '  +4.50%  '
$ws.Range("E51").Value = '  +15.94%  '
